# Atualização 20/07 - Verificação simples da quantidade de batidas inseridas
# antes de importar os dados para o excel.
#
# Adds a new row (row 3) to the "PontoEletrônico" sheet with the next day's
# punch-clock entries: date, entrada, intervalo, retorno intervalo, saída.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "20/07/2023"
$ws.Range("B3").Value = "22:11:37"
$ws.Range("C3").Value = "22:11:38"
$ws.Range("D3").Value = "22:11:39"
$ws.Range("E3").Value = "22:11:40"
